$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra blank row (row 2) - this shifts row3..row7 up by one row,
# so the data previously in A3 moves to A2, and the old A7 disappears.
$ws.Rows("2").Delete()

# Update the active selection to match the post-edit state.
$ws.Range("E6").Select()
